# Action Center Validation Added
# Inserts a new "Processed_PaySlip_Path" setting row into the Tasks/Settings
# sheet (sheet1.xml), pushing the existing rows 14-21 down to rows 15-22,
# and updates the active selection to reflect where the user clicked next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the old row 14 (Output_Excel_Path), shifting
# everything below it down by one row. This also pulls column-A formatting
# down from the row above, matching the target workbook.
$ws.Rows.Item(13).Insert()

# Populate the newly-inserted row with the new configuration entry.
$ws.Range("A13").Value = "Processed_PaySlip_Path"
$ws.Range("B13").Value = "C:\Users\saura\Documents\UiPath\UiPathDemo\Input\Processed\"

# Move the selection to A14, matching the recorded cursor position after
# the edit.
$ws.Range("A14").Select() | Out-Null
